$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 230
$ws.Range("I2").Value = 610
$ws.Range("J2").Value = 2452
$ws.Range("K2").Value = 19
$ws.Range("L2").Value = 678
$ws.Range("M2").Value = 35
$ws.Range("N2").Value = 459
$ws.Range("P2").Value = 7
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 26
$ws.Range("S2").Value = 230
$ws.Range("T2").Value = 455
$ws.Range("U2").Value = 37
$ws.Range("V2").Value = 3734
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 3726
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 53
$ws.Range("AA2").Value = 26
